$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 35 (the empty spacer row between the last data row (34) and the
# totals block that used to start at row 36). This shifts rows 36-38 up to
# 35-37 and Excel automatically adjusts the SUM/formula references.
$ws.Rows.Item(35).Delete()

# Update the active selection to match the post-edit state.
$ws.Range("F34").Select()
